$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value2 = 25000714
$ws.Range("I33").Value2 = 553.1142599999999
$ws.Range("J33").Value2 = 200001840
$ws.Range("K33").Value2 = 553.1142599999999
$ws.Range("L33").Value2 = 200001840
$ws.Range("M33").Value2 = -324.1142599999999
$ws.Range("N33").Value2 = -200002298

# Row 40
$ws.Range("H40").Value2 = 1764.2258
$ws.Range("I40").Value2 = 1557.3158
$ws.Range("J40").Value2 = 2091.8333
$ws.Range("K40").Value2 = 1557.3158
$ws.Range("L40").Value2 = 2091.8333
$ws.Range("M40").Value2 = -1382.3158
$ws.Range("N40").Value2 = -2441.8333

# Row 75
$ws.Range("H75").Value2 = 28950
$ws.Range("I75").Value2 = 28000
$ws.Range("J75").Value2 = 29900
$ws.Range("K75").Value2 = 28000
$ws.Range("L75").Value2 = 29900
$ws.Range("M75").Value2 = -27064
$ws.Range("N75").Value2 = -31772

# Row 78
$ws.Range("H78").Value2 = 28950
$ws.Range("I78").Value2 = 28000
$ws.Range("J78").Value2 = 29900
$ws.Range("K78").Value2 = 84000
$ws.Range("L78").Value2 = 89700
$ws.Range("M78").Value2 = -79320
$ws.Range("N78").Value2 = -99060

# Row 137
$ws.Range("H137").Value2 = 2317438.5
$ws.Range("I137").Value2 = 2700
$ws.Range("J137").Value2 = 3706281.5
$ws.Range("K137").Value2 = 8100
$ws.Range("L137").Value2 = 11118844.5
$ws.Range("M137").Value2 = -5550
$ws.Range("N137").Value2 = -11123944.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("N65").ClearContents()

# Row 86
$ws.Range("H86").Value2 = 1618.8864
$ws.Range("I86").Value2 = 1544.7878
$ws.Range("J86").Value2 = 1841.1818
$ws.Range("K86").Value2 = 1544.7878
$ws.Range("L86").Value2 = 1841.1818
$ws.Range("M86").Value2 = -421.7878000000001
$ws.Range("N86").Value2 = -4087.1818

# Row 89
$ws.Range("H89").Value2 = 1618.8864
$ws.Range("I89").Value2 = 1544.7878
$ws.Range("J89").Value2 = 1841.1818
$ws.Range("K89").Value2 = 7723.939
$ws.Range("L89").Value2 = 9205.909
$ws.Range("M89").Value2 = -2107.939
$ws.Range("N89").Value2 = -20437.909


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 2201.5757
$ws.Range("I31").Value2 = 1383.5
$ws.Range("K31").Value2 = 1383.5
$ws.Range("M31").Value2 = -1088.5

# Row 34
$ws.Range("H34").Value2 = 2201.5757
$ws.Range("I34").Value2 = 1383.5
$ws.Range("K34").Value2 = 1383.5
$ws.Range("M34").Value2 = -1181.5

# Row 122
$ws.Range("H122").Value2 = 1046.5
$ws.Range("I122").Value2 = 857.25
$ws.Range("J122").Value2 = 1425
$ws.Range("K122").Value2 = 2571.75
$ws.Range("L122").Value2 = 4275
$ws.Range("M122").Value2 = -121.75
$ws.Range("N122").Value2 = -9175


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value2 = 1965
$ws.Range("I21").Value2 = 2730.3845
$ws.Range("J21").Value2 = 1633.3334
$ws.Range("K21").Value2 = 8191.1535
$ws.Range("L21").Value2 = 4900.0002
$ws.Range("M21").Value2 = -8018.1535
$ws.Range("N21").Value2 = -5246.0002

# Row 92
$ws.Range("H92").Value2 = 1250660.5
$ws.Range("I92").Value2 = 576
$ws.Range("J92").Value2 = 2500745
$ws.Range("K92").Value2 = 1728
$ws.Range("L92").Value2 = 7502235
$ws.Range("M92").Value2 = -480
$ws.Range("N92").Value2 = -7504731

# Row 94
$ws.Range("H94").Value2 = 4776.1665
$ws.Range("I94").Value2 = 800
$ws.Range("J94").Value2 = 5010.0586
$ws.Range("K94").Value2 = 2400
$ws.Range("L94").Value2 = 15030.1758
$ws.Range("M94").Value2 = -1724
$ws.Range("N94").Value2 = -16382.1758


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value2 = 22378.166
$ws.Range("I4").Value2 = 4168
$ws.Range("J4").Value2 = 40588.332
$ws.Range("K4").Value2 = 4168
$ws.Range("L4").Value2 = 40588.332
$ws.Range("M4").Value2 = -4055
$ws.Range("N4").Value2 = -40814.332

# Row 5
$ws.Range("H5").Value2 = 920
$ws.Range("I5").Value2 = 600
$ws.Range("J5").Value2 = 1000
$ws.Range("K5").Value2 = 600
$ws.Range("L5").Value2 = 1000
$ws.Range("M5").Value2 = -487
$ws.Range("N5").Value2 = -1226

# Row 21
$ws.Range("H21").Value2 = 70007
$ws.Range("J21").Value2 = 70007
$ws.Range("L21").Value2 = 70007
$ws.Range("N21").Value2 = -70355

# Row 24
$ws.Range("H24").Value2 = 1234002.8
$ws.Range("J24").Value2 = 1234002.8
$ws.Range("L24").Value2 = 1234002.8
$ws.Range("N24").Value2 = -1234688.8

# Row 26
$ws.Range("H26").Value2 = 10224
$ws.Range("I26").Value2 = 7995
$ws.Range("J26").Value2 = 11710
$ws.Range("K26").Value2 = 7995
$ws.Range("L26").Value2 = 11710
$ws.Range("M26").Value2 = -7700
$ws.Range("N26").Value2 = -12300

# Row 28
$ws.Range("H28").Value2 = 22378.166
$ws.Range("I28").Value2 = 4168
$ws.Range("J28").Value2 = 40588.332
$ws.Range("K28").Value2 = 4168
$ws.Range("L28").Value2 = 40588.332
$ws.Range("M28").Value2 = -3936
$ws.Range("N28").Value2 = -41052.332

# Row 30
$ws.Range("H30").Value2 = 16849.834
$ws.Range("I30").Value2 = 7999.5
$ws.Range("J30").Value2 = 21275
$ws.Range("K30").Value2 = 7999.5
$ws.Range("L30").Value2 = 21275
$ws.Range("M30").Value2 = -7891.5
$ws.Range("N30").Value2 = -21491

# Row 31
$ws.Range("H31").Value2 = 386.66666
$ws.Range("I31").Value2 = 386.66666
$ws.Range("K31").Value2 = 386.66666
$ws.Range("M31").Value2 = -138.66666

# Row 37
$ws.Range("H37").Value2 = 22378.166
$ws.Range("I37").Value2 = 4168
$ws.Range("J37").Value2 = 40588.332
$ws.Range("K37").Value2 = 4168
$ws.Range("L37").Value2 = 40588.332
$ws.Range("M37").Value2 = -4061
$ws.Range("N37").Value2 = -40802.332

# Row 76
$ws.Range("H76").Value2 = 12857
$ws.Range("J76").Value2 = 12857
$ws.Range("L76").Value2 = 12857
$ws.Range("N76").Value2 = -13533

# Row 79
$ws.Range("H79").Value2 = 12857
$ws.Range("J79").Value2 = 12857
$ws.Range("L79").Value2 = 12857
$ws.Range("N79").Value2 = -15197


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value2 = 4277.5
$ws.Range("I51").Value2 = 4277.5
$ws.Range("K51").Value2 = 4277.5
$ws.Range("M51").Value2 = -3767.5

